$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

# Header label in A1 ("HowMany" becomes the sole shared string).
$ws.Range("A1").Value = "HowMany"

# Values 2..5 in A2:A5.
$ws.Range("A2").Value = 2
$ws.Range("A3").Value = 3
$ws.Range("A4").Value = 4
$ws.Range("A5").Value = 5

# Those cells are unlocked (sheet-protection "input" cells) and carry a
# thin black border on three sides (right/top/bottom only, no left).
$inputRange = $ws.Range("A2:A5")
$inputRange.Locked = $false
$inputRange.Borders.LineStyle = 1
$inputRange.Borders.Color = 0
$inputRange.Borders.Item(7).LineStyle = -4142

# Column A widens to fit the "HowMany" header.
$ws.Columns.Item(1).AutoFit()

# Selection sits on A2; restore the workbook's originally-active sheet
# (Action1) afterwards so the active tab doesn't change.
$ws.Range("A2").Select() | Out-Null
$wb.Worksheets.Item("Action1").Range("A2").Select() | Out-Null
